# Auto-generated Excel COM-interop script
# Applies updated market-price/profit values to the Chocobo Profits workbook sheets
# as produced by the scheduled data-refresh runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 277.4091
$ws.Range("I41").Value = 116
$ws.Range("J41").Value = 324.88235
$ws.Range("K41").Value = 116
$ws.Range("L41").Value = 324.88235
$ws.Range("M41").Value = 324
$ws.Range("N41").Value = -1204.88235

$ws.Range("H98").Value = 6977.857
$ws.Range("I98").Value = 3585.7144
$ws.Range("J98").Value = 10370
$ws.Range("K98").Value = 3585.7144
$ws.Range("L98").Value = 10370
$ws.Range("M98").Value = -2087.7144
$ws.Range("N98").Value = -13366

$ws.Range("H107").Value = 1654.8334
$ws.Range("I107").Value = 2280.5557
$ws.Range("J107").Value = 1029.1111
$ws.Range("K107").Value = 2280.5557
$ws.Range("L107").Value = 1029.1111
$ws.Range("M107").Value = -360.5556999999999
$ws.Range("N107").Value = -4869.1111

$ws.Range("H122").Value = 6977.857
$ws.Range("I122").Value = 3585.7144
$ws.Range("J122").Value = 10370
$ws.Range("K122").Value = 10757.1432
$ws.Range("L122").Value = 31110
$ws.Range("M122").Value = -8307.143199999999
$ws.Range("N122").Value = -36010

$ws.Range("H137").Value = 1589222.5
$ws.Range("I137").Value = 1905754.9
$ws.Range("K137").Value = 5717264.699999999
$ws.Range("M137").Value = -5714714.699999999

$ws.Range("H141").Value = 144538.78
$ws.Range("I141").Value = 183013
$ws.Range("K141").Value = 549039
$ws.Range("M141").Value = -543859

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1431.2222
$ws.Range("I2").Value = 1178.2
$ws.Range("J2").Value = 1747.5
$ws.Range("K2").Value = 1178.2
$ws.Range("L2").Value = 1747.5
$ws.Range("M2").Value = -1065.2
$ws.Range("N2").Value = -1973.5

$ws.Range("H45").Value = 2511.6667
$ws.Range("I45").Value = 2267.5
$ws.Range("K45").Value = 2267.5
$ws.Range("M45").Value = -1890.5

$ws.Range("H61").Value = 1903
$ws.Range("I61").Value = 1302.7778
$ws.Range("K61").Value = 1302.7778
$ws.Range("M61").Value = -1090.7778

$ws.Range("H74").Value = 4729.593
$ws.Range("I74").Value = 6658.643
$ws.Range("K74").Value = 6658.643
$ws.Range("M74").Value = -5784.643

$ws.Range("H77").Value = 4729.593
$ws.Range("I77").Value = 6658.643
$ws.Range("K77").Value = 33293.215
$ws.Range("M77").Value = -28925.215

$ws.Range("H116").Value = 1431.2222
$ws.Range("I116").Value = 1178.2
$ws.Range("J116").Value = 1747.5
$ws.Range("K116").Value = 1178.2
$ws.Range("L116").Value = 1747.5
$ws.Range("M116").Value = 1115.8
$ws.Range("N116").Value = -6335.5

$ws.Range("H132").Value = 2152.077
$ws.Range("I132").Value = 1086.1818
$ws.Range("J132").Value = 2933.7334
$ws.Range("K132").Value = 3258.5454
$ws.Range("L132").Value = 8801.200199999999
$ws.Range("M132").Value = -728.5454
$ws.Range("N132").Value = -13861.2002

$ws.Range("H136").Value = 1903
$ws.Range("I136").Value = 1302.7778
$ws.Range("K136").Value = 3908.3334
$ws.Range("M136").Value = -1358.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1431.2222
$ws.Range("I3").Value = 1178.2
$ws.Range("J3").Value = 1747.5
$ws.Range("K3").Value = 1178.2
$ws.Range("L3").Value = 1747.5
$ws.Range("M3").Value = -1064.2
$ws.Range("N3").Value = -1975.5

$ws.Range("H134").Value = 2401.1853
$ws.Range("I134").Value = 1584.2354
$ws.Range("J134").Value = 3790
$ws.Range("K134").Value = 4752.706200000001
$ws.Range("L134").Value = 11370
$ws.Range("M134").Value = -2217.706200000001
$ws.Range("N134").Value = -16440

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 277.82144
$ws.Range("I7").Value = 255.8
$ws.Range("J7").Value = 303.23077
$ws.Range("K7").Value = 255.8
$ws.Range("L7").Value = 303.23077
$ws.Range("M7").Value = -142.8
$ws.Range("N7").Value = -529.23077

$ws.Range("H31").Value = 3116.8462
$ws.Range("I31").Value = 858.8570999999999
$ws.Range("K31").Value = 858.8570999999999
$ws.Range("M31").Value = -563.8570999999999

$ws.Range("H34").Value = 3116.8462
$ws.Range("I34").Value = 858.8570999999999
$ws.Range("K34").Value = 858.8570999999999
$ws.Range("M34").Value = -656.8570999999999

$ws.Range("H58").Value = 2758.1667
$ws.Range("I58").Value = 1614.7192
$ws.Range("K58").Value = 1614.7192
$ws.Range("M58").Value = -1411.7192

$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 5000
$ws.Range("K62").Value = 5000
$ws.Range("M62").Value = -4376

$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 5000
$ws.Range("K65").Value = 25000
$ws.Range("M65").Value = -21880

$ws.Range("H94").Value = 1645.7858
$ws.Range("I94").Value = 1233.1666
$ws.Range("J94").Value = 1955.25
$ws.Range("K94").Value = 1233.1666
$ws.Range("L94").Value = 1955.25
$ws.Range("M94").Value = -782.1666
$ws.Range("N94").Value = -2857.25

$ws.Range("H132").Value = 3196.7856
$ws.Range("I132").Value = 2794.1667
$ws.Range("J132").Value = 3921.5
$ws.Range("K132").Value = 8382.500100000001
$ws.Range("L132").Value = 11764.5
$ws.Range("M132").Value = -5852.500100000001
$ws.Range("N132").Value = -16824.5

$ws.Range("H134").Value = 2055.1177
$ws.Range("I134").Value = 1161.5
$ws.Range("K134").Value = 3484.5
$ws.Range("M134").Value = -949.5

$ws.Range("H136").Value = 2758.1667
$ws.Range("I136").Value = 1614.7192
$ws.Range("K136").Value = 4844.1576
$ws.Range("M136").Value = -2294.1576

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 25002660
$ws.Range("I80").Value = 41669132
$ws.Range("J80").Value = 2950
$ws.Range("K80").Value = 41669132
$ws.Range("L80").Value = 2950
$ws.Range("M80").Value = -41668134
$ws.Range("N80").Value = -4946

$ws.Range("H83").Value = 25002660
$ws.Range("I83").Value = 41669132
$ws.Range("J83").Value = 2950
$ws.Range("K83").Value = 208345660
$ws.Range("L83").Value = 14750
$ws.Range("M83").Value = -208340668
$ws.Range("N83").Value = -24734

$ws.Range("H122").Value = 4316.737
$ws.Range("I122").Value = 2776.25
$ws.Range("K122").Value = 8328.75
$ws.Range("M122").Value = -5878.75

$ws.Range("H126").Value = 3111.93
$ws.Range("I126").Value = 2828.1375
$ws.Range("J126").Value = 4247.1
$ws.Range("K126").Value = 8484.412499999999
$ws.Range("L126").Value = 12741.3
$ws.Range("M126").Value = -6014.412499999999
$ws.Range("N126").Value = -17681.3

$ws.Range("H132").Value = 4103.5
$ws.Range("I132").Value = 1993.25
$ws.Range("J132").Value = 6213.75
$ws.Range("K132").Value = 5979.75
$ws.Range("L132").Value = 18641.25
$ws.Range("M132").Value = -3449.75
$ws.Range("N132").Value = -23701.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9139.571
$ws.Range("I40").Value = 4004
$ws.Range("J40").Value = 9995.5
$ws.Range("K40").Value = 4004
$ws.Range("L40").Value = 9995.5
$ws.Range("M40").Value = -3868
$ws.Range("N40").Value = -10267.5

$ws.Range("I93").Value = 9260928
$ws.Range("J93").Value = 1708.8334
$ws.Range("K93").Value = 9260928
$ws.Range("L93").Value = 1708.8334
$ws.Range("M93").Value = -9259680
$ws.Range("N93").Value = -4204.8334

$ws.Range("H132").Value = 11013.459
$ws.Range("I132").Value = 11929.583
$ws.Range("K132").Value = 35788.749
$ws.Range("M132").Value = -33258.749

$ws.Range("H136").Value = 4347.522
$ws.Range("I136").Value = 1057.75
$ws.Range("K136").Value = 3173.25
$ws.Range("M136").Value = -623.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6946643
$ws.Range("I132").Value = 1516.7838
$ws.Range("J132").Value = 30307522
$ws.Range("K132").Value = 4550.3514
$ws.Range("L132").Value = 90922566
$ws.Range("M132").Value = -2020.3514
$ws.Range("N132").Value = -90927626

$ws.Range("H136").Value = 4145.12
$ws.Range("I136").Value = 1857.1666
$ws.Range("J136").Value = 10028.429
$ws.Range("K136").Value = 5571.4998
$ws.Range("L136").Value = 30085.287
$ws.Range("M136").Value = -3021.4998
$ws.Range("N136").Value = -35185.287
